# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the existing "Late" / "Heading" / "Outstanding" columns
# one place to the right (N->O, O->P, P->Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate() | Out-Null

# Remember the width of the column immediately to the left (M) so the
# freshly inserted column can inherit it, matching Excel's own behaviour
# when a column is inserted via the UI.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new column at N; everything from N onward shifts right by one.
$ws.Columns("N").Insert() | Out-Null

# Give the new (now empty) column N the inherited width.
$ws.Columns("N").ColumnWidth = $leftWidth

# Restore/update the sheet's active selection.
$ws.Range("T7").Select() | Out-Null
